# Renames the inline-picture shapes (docPr/@name, mirrored onto pic:cNvPr/@name
# by the host) in the document's header/footer stories:
#   - header1.xml  (section 1 "first page" header)  : BTec_Logo-Orange  image2.jpg -> image1.jpg
#   - footer1.xml  (section 1 "first page" footer)  : PearsonLogo       image1.png -> image2.png
#   - footer2.xml  (section 1 "default"    footer)  : PearsonLogo       image1.png -> image2.png

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Header (first page), holds the BTec logo ---
$hdr = $sec.Headers.Item(2)
if ($hdr.Exists) {
    $hdrShapes = $hdr.Range.InlineShapes
    for ($i = 1; $i -le $hdrShapes.Count; $i++) {
        $shape = $hdrShapes.Item($i)
        if ($shape.AlternativeText -eq "BTec_Logo-Orange") {
            $shape.Name = "image1.jpg"
        }
    }
}

# --- Footers (default + first page), hold the Pearson Edexcel logo ---
for ($fi = 1; $fi -le 3; $fi++) {
    $ftr = $sec.Footers.Item($fi)
    if ($ftr.Exists) {
        $ftrShapes = $ftr.Range.InlineShapes
        for ($i = 1; $i -le $ftrShapes.Count; $i++) {
            $shape = $ftrShapes.Item($i)
            if ($shape.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shape.Name = "image2.png"
            }
        }
    }
}
